$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.532.81"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.840.98"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'313.62"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.4253"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.3671"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "'0.07278"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'0.8711"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").Value = "'20.83"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "1.831.64"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "'5.404"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'6.531"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'0.06942"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'80.45"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "'0.000009029"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'15.50"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "27.592.24"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'5.066"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "'10.87"
$ws.Range("E23").Value = "  +4.91%  "
$ws.Range("D24").Value = "2.073.73"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "'1.956"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").Value = "'154.28"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "'5.259"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'115.43"
$ws.Range("E29").Value = "  -5.60%  "
$ws.Range("D30").Value = "'1.861"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'0.08876"
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'0.7768"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'4.551"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'2.957"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'1.155"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "'1.102"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").Value = "'0.05396"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "'0.01948"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "'2.817"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'0.5129"
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("D42").Value = "'0.1667"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "'6.747"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").Value = "'8.553"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "'10.56"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "'106.70"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "'0.06536"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'0.4726"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "'1.641"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "'1.797"
$ws.Range("E51").Value = "  +4.44%  "
